$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 8 (diff @ 1018)
$ws.Range("H8").Value = 1345.25
$ws.Range("I8").Value = 49.57143
$ws.Range("J8").Value = 3159.2
$ws.Range("K8").Value = 148.71429
$ws.Range("L8").Value = 9477.599999999999
$ws.Range("M8").Value = -9.714290000000005
$ws.Range("N8").Value = -9755.599999999999

$ws = $wb.Worksheets.Item("ARM")
# row 44 (diff @ 9952)
$ws.Range("H44").Value = 30342.334
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 30342.334
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 30342.334
$ws.Range("N44").Value = -31318.334

# row 55 (diff @ 10485)
$ws.Range("H55").Value = 23655.334
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 23655.334
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 23655.334
$ws.Range("N55").Value = -24285.334

# row 80 (diff @ 11695)
$ws.Range("H80").Value = 34097.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 34097.5
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 34097.5
$ws.Range("N80").Value = -36093.5

# row 83 (diff @ 11839)
$ws.Range("H83").Value = 34097.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 34097.5
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 102292.5
$ws.Range("N83").Value = -112276.5

# row 122 (diff @ 13753)
$ws.Range("H122").Value = 1541.5834
$ws.Range("I122").Value = 1496
$ws.Range("J122").Value = 1605.4
$ws.Range("K122").Value = 4488
$ws.Range("L122").Value = 4816.200000000001
$ws.Range("M122").Value = -2038
$ws.Range("N122").Value = -9716.200000000001

$ws = $wb.Worksheets.Item("BSM")
# row 35 (diff @ 16456)
$ws.Range("H35").Value = 29916.285
$ws.Range("I35").Value = 15000
$ws.Range("J35").Value = 32402.334
$ws.Range("K35").Value = 15000
$ws.Range("L35").Value = 32402.334
$ws.Range("M35").Value = -14690
$ws.Range("N35").Value = -33022.334

# row 56 (diff @ 17485)
$ws.Range("H56").Value = 30110
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 30110
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 30110
$ws.Range("N56").Value = -31588

# row 82 (diff @ 18765)
$ws.Range("H82").Value = 17002.785
$ws.Range("I82").Value = 3333.3333
$ws.Range("J82").Value = 27254.875
$ws.Range("K82").Value = 3333.3333
$ws.Range("L82").Value = 27254.875
$ws.Range("M82").Value = -2950.3333
$ws.Range("N82").Value = -28020.875

# row 85 (diff @ 18918)
$ws.Range("H85").Value = 17002.785
$ws.Range("I85").Value = 3333.3333
$ws.Range("J85").Value = 27254.875
$ws.Range("K85").Value = 3333.3333
$ws.Range("L85").Value = 27254.875
$ws.Range("M85").Value = -2007.3333
$ws.Range("N85").Value = -29906.875

$ws = $wb.Worksheets.Item("CRP")
# row 16 (diff @ 22515)
$ws.Range("H16").Value = 1739
$ws.Range("I16").Value = 1678.2
$ws.Range("J16").Value = 1799.8
$ws.Range("K16").Value = 1678.2
$ws.Range("L16").Value = 1799.8
$ws.Range("M16").Value = -1391.2
$ws.Range("N16").Value = -2373.8

# row 41 (diff @ 23767)
$ws.Range("H41").Value = 22030.6
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 22030.6
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 22030.6
$ws.Range("N41").Value = -22886.6
$ws.Range("M41").ClearContents()

# row 50 (diff @ 24217)
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

# row 51 (diff @ 24266)
$ws.Range("H51").Value = 29574.8
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 29574.8
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 29574.8
$ws.Range("N51").Value = -31046.8
$ws.Range("M51").ClearContents()

# row 58 (diff @ 24612)
$ws.Range("H58").Value = 909.7308
$ws.Range("I58").Value = 988.75
$ws.Range("J58").Value = 646.3333
$ws.Range("K58").Value = 988.75
$ws.Range("L58").Value = 646.3333
$ws.Range("M58").Value = -785.75
$ws.Range("N58").Value = -1052.3333

# row 59 (diff @ 24664)
$ws.Range("H59").Value = 333353180
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 333353180
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 333353180
$ws.Range("N59").Value = -333355470

# row 60 (diff @ 24713)
$ws.Range("H60").Value = 12917.6
$ws.Range("I60").Value = 5771.5
$ws.Range("J60").Value = 17681.666
$ws.Range("K60").Value = 5771.5
$ws.Range("L60").Value = 17681.666
$ws.Range("M60").Value = -5260.5
$ws.Range("N60").Value = -18703.666

# row 61 (diff @ 24765)
$ws.Range("H61").Value = 29574.8
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 29574.8
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 29574.8
$ws.Range("N61").Value = -30270.8
$ws.Range("M61").ClearContents()

# row 68 (diff @ 25117)
$ws.Range("H68").Value = 25880.715
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 25880.715
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 25880.715
$ws.Range("N68").Value = -27378.715

# row 71 (diff @ 25267)
$ws.Range("H71").Value = 25880.715
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 25880.715
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 77642.145
$ws.Range("N71").Value = -85130.145

# row 109 (diff @ 27141)
$ws.Range("H109").Value = 23052.285
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 23052.285
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 23052.285
$ws.Range("N109").Value = -25132.285

# row 113 (diff @ 27334)
$ws.Range("H113").Value = 1739
$ws.Range("I113").Value = 1678.2
$ws.Range("J113").Value = 1799.8
$ws.Range("K113").Value = 1678.2
$ws.Range("L113").Value = 1799.8
$ws.Range("M113").Value = 491.8
$ws.Range("N113").Value = -6139.8

# row 136 (diff @ 28470)
$ws.Range("H136").Value = 909.7308
$ws.Range("I136").Value = 988.75
$ws.Range("J136").Value = 646.3333
$ws.Range("K136").Value = 2966.25
$ws.Range("L136").Value = 1938.9999
$ws.Range("M136").Value = -416.25
$ws.Range("N136").Value = -7038.9999

$ws = $wb.Worksheets.Item("CUL")
# row 32 (diff @ 30382)
$ws.Range("H32").Value = 3500
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 3500
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 10500
$ws.Range("N32").Value = -11066

# row 46 (diff @ 31089)
$ws.Range("H46").Value = 2155
$ws.Range("I46").Value = 280
$ws.Range("J46").Value = 3196.6667
$ws.Range("K46").Value = 840
$ws.Range("L46").Value = 9590.000100000001
$ws.Range("M46").Value = -749

# row 64 (diff @ 32004)
$ws.Range("H64").Value = 19336.334
$ws.Range("I64").Value = 2000
$ws.Range("J64").Value = 28004.5
$ws.Range("K64").Value = 6000
$ws.Range("L64").Value = 84013.5
$ws.Range("M64").Value = -5730
$ws.Range("N64").Value = -84553.5

# row 67 (diff @ 32157)
$ws.Range("H67").Value = 19336.334
$ws.Range("I67").Value = 2000
$ws.Range("J67").Value = 28004.5
$ws.Range("K67").Value = 6000
$ws.Range("L67").Value = 84013.5
$ws.Range("M67").Value = -5064
$ws.Range("N67").Value = -85885.5

# row 68 (diff @ 32209)
$ws.Range("H68").Value = 1536.2059
$ws.Range("I68").Value = 983.7857
$ws.Range("J68").Value = 1922.9
$ws.Range("K68").Value = 2951.3571
$ws.Range("L68").Value = 5768.700000000001
$ws.Range("M68").Value = -2140.3571
$ws.Range("N68").Value = -7390.700000000001

# row 69 (diff @ 32261)
$ws.Range("H69").Value = 4989.5
$ws.Range("I69").Value = 500
$ws.Range("J69").Value = 5488.3335
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 16465.0005
$ws.Range("M69").Value = -689
$ws.Range("N69").Value = -18087.0005

# row 71 (diff @ 32365)
$ws.Range("H71").Value = 1536.2059
$ws.Range("I71").Value = 983.7857
$ws.Range("J71").Value = 1922.9
$ws.Range("K71").Value = 8854.0713
$ws.Range("L71").Value = 17306.1
$ws.Range("M71").Value = -4798.0713
$ws.Range("N71").Value = -25418.1

# row 72 (diff @ 32417)
$ws.Range("H72").Value = 4989.5
$ws.Range("I72").Value = 500
$ws.Range("J72").Value = 5488.3335
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 49395.0015
$ws.Range("M72").Value = -444
$ws.Range("N72").Value = -57507.0015

# row 81 (diff @ 32885)
$ws.Range("H81").Value = 250001020
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 250001020
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 750003060
$ws.Range("N81").Value = -750005306
$ws.Range("M81").ClearContents()

# row 84 (diff @ 33038)
$ws.Range("H84").Value = 250001020
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 250001020
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 2250009180
$ws.Range("N84").Value = -2250020412
$ws.Range("M84").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# row 3 (diff @ 36131)
$ws.Range("H3").Value = 687.25
$ws.Range("I3").Value = 98.5
$ws.Range("J3").Value = 1276
$ws.Range("K3").Value = 98.5
$ws.Range("L3").Value = 1276
$ws.Range("M3").Value = 17.5
$ws.Range("N3").Value = -1508

# row 57 (diff @ 38771)
$ws.Range("H57").Value = 21477.924
$ws.Range("I57").Value = 3000
$ws.Range("J57").Value = 23017.75
$ws.Range("K57").Value = 3000
$ws.Range("L57").Value = 23017.75
$ws.Range("M57").Value = -2180
$ws.Range("N57").Value = -24657.75

$ws = $wb.Worksheets.Item("LTW")
# row 43 (diff @ 45006)
$ws.Range("H43").Value = 4000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 4000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386

# row 109 (diff @ 48228)
$ws.Range("H109").Value = 30037.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 30037.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 30037.5
$ws.Range("N109").Value = -32811.5

# row 136 (diff @ 49560)
$ws.Range("H136").Value = 1508.8148
$ws.Range("I136").Value = 1350.1305
$ws.Range("J136").Value = 2421.25
$ws.Range("K136").Value = 4050.3915
$ws.Range("L136").Value = 7263.75
$ws.Range("M136").Value = -1500.3915
$ws.Range("N136").Value = -12363.75

$ws = $wb.Worksheets.Item("WVR")
# row 41 (diff @ 51859)
$ws.Range("H41").Value = 16610.8
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 16610.8
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 16610.8
$ws.Range("N41").Value = -17390.8

# row 109 (diff @ 55173)
$ws.Range("H109").Value = 31050.8
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 31050.8
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 31050.8
$ws.Range("N109").Value = -33824.8
